$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string / label edits -------------------------------------------------
# Status value "ffxqa" -> "open" (appears twice, in the two "Status" mini-tables)
$ws.Range("A16").Value = "open"
$ws.Range("A25").Value = "open"

# Column headers: insert "Customer Generated Ticket" ahead of "Customer Testing",
# and drop "FFX Testing" (both header rows of the "CAC/MOF Requestor" tables)
$ws.Range("C7").Value  = "Customer Generated Ticket"
$ws.Range("D7").Value  = "Customer Testing"
$ws.Range("C24").Value = "Customer Generated Ticket"
$ws.Range("D24").Value = "Customer Testing"

# --- Numeric value edits ---------------------------------------------------------
$ws.Range("B2").Value  = 1
$ws.Range("B4").Value  = 1

$ws.Range("C8").Value  = 1
$ws.Range("D8").Value  = 1
$ws.Range("F8").Value  = 6

$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1
$ws.Range("F10").Value = 6

$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 3
$ws.Range("F11").Value = 13

$ws.Range("B16").Value = 0
$ws.Range("D16").Value = 1

$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 3

$ws.Range("B20").Value = 6
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 13

$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 0

$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 3
$ws.Range("F26").Value = 3

$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 3
$ws.Range("F29").Value = 13

$ws.Range("B35").Value = 2
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 3

$ws.Range("B36").Value = 2
$ws.Range("D36").Value = 1

$ws.Range("C37").Value = 2
$ws.Range("E37").Value = 5

$ws.Range("B38").Value = 6
$ws.Range("C38").Value = 4
$ws.Range("E38").Value = 13

$ws.Range("B43").Value = 2
$ws.Range("D43").Value = 3
$ws.Range("E43").Value = 6

$ws.Range("B44").Value = 2
$ws.Range("D44").Value = 2
$ws.Range("E44").Value = 4

$ws.Range("B46").Value = 6
$ws.Range("D46").Value = 6
$ws.Range("E46").Value = 13

$ws.Range("B51").Value = 8
$ws.Range("B52").Value = 13

$ws.Range("B57").Value = 2
$ws.Range("B58").Value = 2
$ws.Range("B60").Value = 6

$ws.Range("B64").Value = 3
$ws.Range("B66").Value = 5
$ws.Range("B67").Value = 13
